$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows right before row 205; this pushes the existing
# rows 205-210 down to 208-213 (with all their data/formatting intact).
$ws.Rows("205:207").Insert()

# Row 205: new weekly entry
$ws.Cells.Item(205, 1).Value = 4
$ws.Cells.Item(205, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(205, 3).Value = "Los Lagos"
$ws.Cells.Item(205, 4).Value = 44585
$ws.Cells.Item(205, 5).Value = 10
$ws.Cells.Item(205, 6).Value = 100112028
$ws.Cells.Item(205, 7).Value = "Sandia"
$ws.Cells.Item(205, 8).Value = "Sin especificar"
$ws.Cells.Item(205, 9).Value = "Primera"
$ws.Cells.Item(205, 10).Value = 3000
$ws.Cells.Item(205, 11).Value = 3000
$ws.Cells.Item(205, 12).Value = 3000
$ws.Cells.Item(205, 13).Value = 3000
$ws.Cells.Item(205, 14).Value = '$/unidad'
$ws.Cells.Item(205, 15).Value = "Región del Maule"
$ws.Cells.Item(205, 16).Value = 3000
$ws.Cells.Item(205, 17).Value = 1
$ws.Cells.Item(205, 18).Value = "Hortaliza"

# Row 206: new weekly entry
$ws.Cells.Item(206, 1).Value = 4
$ws.Cells.Item(206, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(206, 3).Value = "Los Lagos"
$ws.Cells.Item(206, 4).Value = 44585
$ws.Cells.Item(206, 5).Value = 10
$ws.Cells.Item(206, 6).Value = 100112028
$ws.Cells.Item(206, 7).Value = "Sandia"
$ws.Cells.Item(206, 8).Value = "Sin especificar"
$ws.Cells.Item(206, 9).Value = "Segunda"
$ws.Cells.Item(206, 10).Value = 3000
$ws.Cells.Item(206, 11).Value = 2500
$ws.Cells.Item(206, 12).Value = 2500
$ws.Cells.Item(206, 13).Value = 2500
$ws.Cells.Item(206, 14).Value = '$/unidad'
$ws.Cells.Item(206, 15).Value = "Región del Maule"
$ws.Cells.Item(206, 16).Value = 2500
$ws.Cells.Item(206, 17).Value = 1
$ws.Cells.Item(206, 18).Value = "Hortaliza"

# Row 207: new weekly entry
$ws.Cells.Item(207, 1).Value = 4
$ws.Cells.Item(207, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(207, 3).Value = "Los Lagos"
$ws.Cells.Item(207, 4).Value = 44585
$ws.Cells.Item(207, 5).Value = 10
$ws.Cells.Item(207, 6).Value = 100112028
$ws.Cells.Item(207, 7).Value = "Sandia"
$ws.Cells.Item(207, 8).Value = "Sin especificar"
$ws.Cells.Item(207, 9).Value = "Tercera"
$ws.Cells.Item(207, 10).Value = 3000
$ws.Cells.Item(207, 11).Value = 2000
$ws.Cells.Item(207, 12).Value = 2000
$ws.Cells.Item(207, 13).Value = 2000
$ws.Cells.Item(207, 14).Value = '$/unidad'
$ws.Cells.Item(207, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(207, 16).Value = 2000
$ws.Cells.Item(207, 17).Value = 1
$ws.Cells.Item(207, 18).Value = "Hortaliza"
